# Insert a new weekly price row above row 592 (pushing the existing
# rows 592:621 down to 593:622) and populate the new row with the
# latest week's data for Cebollín @ Femacal de La Calera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 592:621 down to 593:622, leaving a blank row 592.
$ws.Rows("592:592").Insert()

# Populate the new row 592 with this week's data.
$ws.Range("A592").Value = 3
$ws.Range("B592").Value = "Femacal de La Calera"
$ws.Range("C592").Value = "Coquimbo"
$ws.Range("D592").Value = 44939
$ws.Range("E592").Value = 5
$ws.Range("F592").Value = 100112037
$ws.Range("G592").Value = "Cebollín"
$ws.Range("H592").Value = "Sin especificar"
$ws.Range("I592").Value = "Primera"
$ws.Range("J592").Value = 230
$ws.Range("K592").Value = 3000
$ws.Range("L592").Value = 3300
$ws.Range("M592").Value = 3143
$ws.Range("N592").Value = "$/paquete 36 unidades"
$ws.Range("O592").Value = "Provincia de Quillota"
$ws.Range("P592").Value = 87
$ws.Range("Q592").Value = 36
$ws.Range("R592").Value = "Hortaliza"
